# Add a new roster entry (Cheong Ming Lun / m.cheong@digipen.edu) to row 7,
# reusing row 4's formatting (including its hyperlink) so the new cells pick
# up the existing "Name"/"E-mail" styles instead of minting new ones, then
# re-point the copied hyperlink at the new person's mailto address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (and hyperlink) from row 4 into row 7.
$ws.Range("A4:B4").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)

# Fill in the new team member's details.
$ws.Range("A7").Value = "Cheong Ming Lun"
$ws.Range("B7").Value = "m.cheong@digipen.edu"

# Re-target the hyperlink that was copied along with B4's formatting.
$ws.Range("B7").Hyperlinks(1).Address = "mailto:m.cheong@digipen.edu"

# Match the last-saved selection recorded in the workbook.
[void]$ws.Range("I10").Select()
